# Added New Mac-Address and Document Types
# Append a new data row (row 33) to the master-reg_center_machine sheet,
# reusing the existing lang_code/is_active/cr_by/cr_dtimes values, and
# update the sheet's viewport/selection like the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Match the updated viewport/selection recorded after the edit.
$ws.Range("E29").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
